$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# Target stored width (OOXML <col width>) is 13.4101848602295 "characters".
# Excel's COM ColumnWidth setter quantizes to 1/6-character pixel steps
# (stored = Round(input*6 + 5) / 6), so the nearest reachable grid value is
# 80/6 = 13.333333333333334 -- reached by feeding it 12.5 (the exact middle
# of the input bucket that rounds to that grid point).
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
